# Checkpoint from VS Code for cloud agent session
# Adds reminder-tracking columns (N:P) to the Putt Allotment sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells N1:P1 -------------------------------------------------
# Copy the formatting from the existing last header cell (M1 - "STATUS", style
# index 1: bold, bordered, centered) onto the three new header cells, then set
# their text.
$ws.Range("M1").Copy()
$ws.Range("N1:P1").PasteSpecial(-4122)

$ws.Range("N1").Value = "REMINDER_ROW_ID"
$ws.Range("O1").Value = "REMINDER_SNOOZE_UNTIL"
$ws.Range("P1").Value = "REMINDER_DISMISSED"

# --- Per-row reminder metadata ----------------------------------------------
# Column N: stable UUID identifying each data row for the reminder feature.
# Column O: snooze-until timestamp - left blank for every existing row.
# Column P: dismissed flag - False for every existing row.

$rowIds = @{
    2 = "11b46c4e-084f-42b1-acd9-7a800fb62cf4"
    3 = "d5caa7d8-6422-4365-ba4c-6033a58cf91a"
    4 = "4b013cf0-75ea-4221-bc28-0d683724c08a"
    5 = "ba135f90-2877-472c-b0af-fffbab774325"
    6 = "77bb20c2-2f4a-43fe-9275-a10bde9b17f9"
}

foreach ($r in 2..6) {
    # Column O starts out blank - copy formatting/blankness from an existing
    # empty (unstyled) data cell so the cell is materialised rather than left
    # absent.
    $ws.Cells.Item($r, 9).Copy()
    $ws.Cells.Item($r, 15).PasteSpecial(-4122)

    $ws.Cells.Item($r, 14).Value = $rowIds[$r]
    $ws.Cells.Item($r, 16).Value = $false
}
